$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$helper.Value = "2.85768"
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$helper.Value = "0.986547"
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$helper.Value = "1.8515"
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$helper.Value = "1.22979"
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$helper.Value = "1.04718"
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$helper.Value = "1.88103"
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$helper.Value = "0.825483"
$helper.Copy()
$ws.Range("E2").PasteSpecial(-4163)
$helper.Value = "2.43073"
$helper.Copy()
$ws.Range("E3").PasteSpecial(-4163)
$helper.Value = "0.152433"
$helper.Copy()
$ws.Range("E4").PasteSpecial(-4163)
$helper.Value = "2.06232"
$helper.Copy()
$ws.Range("E5").PasteSpecial(-4163)
$helper.Value = "2.35903"
$helper.Copy()
$ws.Range("E6").PasteSpecial(-4163)
$helper.Value = "1.02378"
$helper.Copy()
$ws.Range("E7").PasteSpecial(-4163)
$helper.Value = "2.20988"
$helper.Copy()
$ws.Range("E8").PasteSpecial(-4163)
$helper.Value = "1.83664"
$helper.Copy()
$ws.Range("E9").PasteSpecial(-4163)
$helper.Value = "0.602734"
$helper.Copy()
$ws.Range("E10").PasteSpecial(-4163)
$helper.Value = "2.66831"
$helper.Copy()
$ws.Range("E11").PasteSpecial(-4163)
$helper.Value = "0.603434"
$helper.Copy()
$ws.Range("E12").PasteSpecial(-4163)
$helper.Value = "2.05417"
$helper.Copy()
$ws.Range("E13").PasteSpecial(-4163)
$helper.Value = "1.41457"
$helper.Copy()
$ws.Range("E14").PasteSpecial(-4163)
$helper.Value = "1.45797"
$helper.Copy()
$ws.Range("E15").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("F2").PasteSpecial(-4163)
$helper.Value = "0"
$helper.Copy()
$ws.Range("G2").PasteSpecial(-4163)
$helper.Value = "0"
$helper.Copy()
$ws.Range("F3").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("G3").PasteSpecial(-4163)
$helper.Value = "3"
$helper.Copy()
$ws.Range("F4").PasteSpecial(-4163)
$helper.Value = "0"
$helper.Copy()
$ws.Range("G4").PasteSpecial(-4163)
$helper.Value = "1"
$helper.Copy()
$ws.Range("F5").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("G5").PasteSpecial(-4163)
$helper.Value = "3"
$helper.Copy()
$ws.Range("F6").PasteSpecial(-4163)
$helper.Value = "3"
$helper.Copy()
$ws.Range("G6").PasteSpecial(-4163)
$helper.Value = "1"
$helper.Copy()
$ws.Range("F7").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("G7").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("F8").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("G8").PasteSpecial(-4163)
$helper.Value = "1"
$helper.Copy()
$ws.Range("F9").PasteSpecial(-4163)
$helper.Value = "3"
$helper.Copy()
$ws.Range("G9").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("F10").PasteSpecial(-4163)
$helper.Value = "0"
$helper.Copy()
$ws.Range("G10").PasteSpecial(-4163)
$helper.Value = "1"
$helper.Copy()
$ws.Range("F11").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("G11").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("F12").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("G12").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("F13").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("G13").PasteSpecial(-4163)
$helper.Value = "1"
$helper.Copy()
$ws.Range("F14").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("G14").PasteSpecial(-4163)
$helper.Value = "3"
$helper.Copy()
$ws.Range("F15").PasteSpecial(-4163)
$helper.Value = "2"
$helper.Copy()
$ws.Range("G15").PasteSpecial(-4163)

$helper.Clear()
$excel.CutCopyMode = $false
